$d = $word.ActiveDocument

# Helper: force a hard run boundary around $rng by toggling a formatting
# property on and back off. The engine coalesces freshly-edited adjacent
# runs that share identical formatting, so an explicit (even if net-zero)
# formatting operation is needed to keep the edited text as its own run.
function Split-At($rng) {
    $rng.Bold = 1
    $rng.Bold = 0
}

# ---------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that originally sat in the title
#    paragraph ("... Penetration of Neutrons through shielding").
# ---------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------
# 2. "April 2020" -> "May 2020" on the title-page date line only (leave
#    the "Accessed 20 April. 2020" citation further down untouched).
# ---------------------------------------------------------------------
$dateRange = $d.Paragraphs(10).Range
Write-Host "Date paragraph before: " $dateRange.Text
$dateRange.Find.Execute("April ", $true, $false, $false, $false, $false, $true, 1, $false, "May ", 2) | Out-Null

# ---------------------------------------------------------------------
# 3. Expand the Abstract paragraph with the additional sentences,
#    keeping "...a neutron" and the new continuation as two runs.
# ---------------------------------------------------------------------
$abs = $d.Content
$abs.Find.Execute("The analysis of the behaviour of a neutron ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$abs.Text = "The analysis of the behaviour of a neutron"

$absTail = $abs.Duplicate
$absTail.Collapse(0)
$absTail.InsertAfter(" penetrating through shielding has been conducted. This was done by modelling the neutron as a random walk with the utilization of Monte Carlo methods. The materials that were tested included water, lead and graphite.")

$absTailFound = $d.Content
$absTailFound.Find.Execute(" penetrating through shielding has been conducted. This was done by modelling the neutron as a random walk with the utilization of Monte Carlo methods. The materials that were tested included water, lead and graphite.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Split-At $absTailFound

# ---------------------------------------------------------------------
# 4. Re-add a "_GoBack" bookmark at the point where the author's cursor
#    was left - right after "...damping coefficient." in section 4.2.
# ---------------------------------------------------------------------
$tailSentence = $d.Content
$tailSentence.Find.Execute("damping forces present in the system arising from the damping coefficient.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$tailSentence.InsertAfter("~GOBACKMARK~")

$markerRange = $d.Content
$markerRange.Find.Execute("~GOBACKMARK~", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $markerRange)
$markerRange.Text = ""

# ---------------------------------------------------------------------
# 5. Reference list: split a handful of names/words into their own runs
#    (mirrors the spell-checker "Süli"/"Mayers"/etc. spans in the diff).
# ---------------------------------------------------------------------

# [5] Suli, E. and Mayers, D., 2014.
$ref5 = $d.Content
$ref5.Find.Execute("[5] Süli, E. and Mayers, D., 2014. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s = $ref5.Start
Split-At ($d.Range($s + 4, $s + 8))
Split-At ($d.Range($s + 17, $s + 23))

# [7] Press, W. and Vetterling, W., 2007.
$ref7 = $d.Content
$ref7.Find.Execute("[7] Press, W. and Vetterling, W., 2007. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s = $ref7.Start
Split-At ($d.Range($s + 18, $s + 28))

# [8] Verlet, L., 1967. Computer "Experiments" ... Thermodynamical Properties ...
$ref8 = $d.Content
$ref8.Find.Execute('[8] Verlet, L., 1967. Computer "Experiments" on Classical Fluids. I. Thermodynamical Properties of Lennard-Jones Molecules. ', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s = $ref8.Start
Split-At ($d.Range($s + 69, $s + 84))

# [9] Iott, J., Haftka, R. and Adelman, H., 1985.
$ref9 = $d.Content
$ref9.Find.Execute("[9] Iott, J., Haftka, R. and Adelman, H., 1985. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s = $ref9.Start
Split-At ($d.Range($s + 4, $s + 8))
Split-At ($d.Range($s + 14, $s + 20))

# [10] ... SplitSky/Scientific_Programming (italic title of the repo)
$ref10 = $d.Content
$ref10.Find.Execute("SplitSky/Scientific_Programming", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s = $ref10.Start
Split-At ($d.Range($s, $s + 8))
Split-At ($d.Range($s + 9, $s + 31))

Write-Host "edit complete"
